# Adds the "line" (table 9) and the start of the "linespacing" (table 10)
# sections to the "Significado dos índices" sheet, mirroring the formatting
# of the existing tables above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Copy formatting for the new rows -------------------------------------
# Rows 32-46 (table 9 "trans/line/insert") reuse the style of the
# "bustype" table rows (style index 1 / orange banding).
$ws.Range("A22:H22").Copy()
$ws.Range("A32:H46").PasteSpecial(-4122)

# Row 47 (first row of table 10 "trans/linespacing/insert") reuses the
# style of the "busbustyperel" table rows (style index 2 / blue banding).
$ws.Range("A27:H27").Copy()
$ws.Range("A47:H47").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column B: "Significado" for the new "line" rows -----------------------
$ws.Cells.Item(32, 2).Value = "caseID"
$ws.Cells.Item(33, 2).Value = "inicialBus"
$ws.Cells.Item(34, 2).Value = "finalBus"
$ws.Cells.Item(35, 2).Value = "sequencialNumber"
$ws.Cells.Item(36, 2).Value = "resistence"
$ws.Cells.Item(37, 2).Value = "reactance"
$ws.Cells.Item(38, 2).Value = "susceptance"
$ws.Cells.Item(39, 2).Value = "MVAratings#1"
$ws.Cells.Item(40, 2).Value = "MVAratings#2"
$ws.Cells.Item(41, 2).Value = "MVAratings#3"
$ws.Cells.Item(42, 2).Value = "Description"
$ws.Cells.Item(43, 2).Value = "circuitNumber"
$ws.Cells.Item(44, 2).Value = "area"

# --- Column A: "Label" for the new "line" rows ------------------------------
$ws.Cells.Item(32, 1).Value = "label30"
$ws.Cells.Item(33, 1).Value = "label31"
$ws.Cells.Item(34, 1).Value = "label32"
$ws.Cells.Item(35, 1).Value = "label33"
$ws.Cells.Item(36, 1).Value = "label34"
$ws.Cells.Item(37, 1).Value = "label39"
$ws.Cells.Item(38, 1).Value = "label38"
$ws.Cells.Item(39, 1).Value = "label37"
$ws.Cells.Item(40, 1).Value = "label36"
$ws.Cells.Item(41, 1).Value = "label35"
$ws.Cells.Item(42, 1).Value = "label40"
$ws.Cells.Item(43, 1).Value = "label42"
$ws.Cells.Item(44, 1).Value = "label41"

# --- Column C: "De qual painel?" (comboboxes first, then textboxes) --------
$ws.Cells.Item(32, 3).Value = "combobox9"
$ws.Cells.Item(44, 3).Value = "combobox10"
$ws.Cells.Item(35, 3).Value = "texbox15"
$ws.Cells.Item(36, 3).Value = "texbox16"
$ws.Cells.Item(37, 3).Value = "texbox17"
$ws.Cells.Item(38, 3).Value = "texbox18"
$ws.Cells.Item(39, 3).Value = "texbox19"
$ws.Cells.Item(40, 3).Value = "texbox20"
$ws.Cells.Item(41, 3).Value = "texbox21"
$ws.Cells.Item(42, 3).Value = "texbox22"
$ws.Cells.Item(43, 3).Value = "texbox23"

# --- Column D/E: table id and "Significado do painel" for all rows ---------
for ($r = 32; $r -le 46; $r++) {
    $ws.Cells.Item($r, 4).Value = 9
    $ws.Cells.Item($r, 5).Value = "trans/line/insert"
}

# --- Row 45: Submit button --------------------------------------------------
$ws.Cells.Item(45, 1).Value = "button10"
$ws.Cells.Item(45, 2).Value = "Submit"
$ws.Cells.Item(45, 3).Value = "-"

# --- Row 46: Clear button ---------------------------------------------------
$ws.Cells.Item(46, 1).Value = "button11"
$ws.Cells.Item(46, 2).Value = "clear"
$ws.Cells.Item(46, 3).Value = "-"

# --- Fill in the remaining comboboxes for rows 33-34 ------------------------
$ws.Cells.Item(33, 3).Value = "combobox11"
$ws.Cells.Item(34, 3).Value = "combobox12"

# --- Table 10: trans/linespacing/insert (first row only) -------------------
$ws.Cells.Item(47, 1).Value = "label43"
$ws.Cells.Item(47, 2).Value = "FormNotUsed "
$ws.Cells.Item(47, 4).Value = 10
$ws.Cells.Item(47, 5).Value = "trans/linespacing/insert"

# --- Update the view so the new rows are visible ----------------------------
$ws.Range("E48").Select()
$excel.ActiveWindow.ScrollRow = 30
